$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the contents of the data rows (A2:N3) - simulates selecting the
# range and pressing Delete: values are removed, cell formatting (styles)
# is preserved.
$range = $ws.Range("A2:N3")
$range.ClearContents()

# Update the active selection to match the post-edit state (A2 active,
# A2:N3 selected).
$ws.Range("A2:N3").Select()
